$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row (option A/B/C/D, correct option -> snake_case names) ---
$ws.Range("C1").Value = "option_a"
$ws.Range("D1").Value = "option_b"
$ws.Range("E1").Value = "option_c"
$ws.Range("F1").Value = "option_d"
$ws.Range("G1").Value = "correct_option"

# --- Remove the bulk of the quiz questions, keeping only the first two data rows ---
# Clear the content of rows 4 through 11 (all columns), leaving the rows themselves
$ws.Range("A4:H11").ClearContents()

# Remove the old trailing placeholder row (previously row 12)
$ws.Rows("12").Delete()

# --- Re-apply a (new) custom number format to the category_id column ---
$ws.Range("A2:A11").NumberFormat = "0;[Red]0"

# --- Adjust column widths to better fit the now much smaller data set ---
$ws.Columns("A").ColumnWidth = 11.333333333333332
$ws.Columns("B").ColumnWidth = 57
$ws.Columns("F").ColumnWidth = 13.833333333333334
$ws.Columns("G").ColumnWidth = 13.666666666666666

# --- Update the selected/active cell ---
$null = $ws.Range("A2").Select()

# --- Set the page to print in portrait orientation ---
$ws.PageSetup.Orientation = 1

Write-Output "done"
